$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 175, shifting rows 175..253 down to 176..254.
$ws.Rows.Item(175).EntireRow.Insert()

# Populate the newly inserted row 175 with its data.
$ws.Range("A175").Value = 11
$ws.Range("B175").Value = "Vega Monumental Concepción"
$ws.Range("C175").Value = "Bíobío"
$ws.Range("D175").Value = 44582
$ws.Range("E175").Value = 8
$ws.Range("F175").Value = "Fruta"
$ws.Range("G175").Value = 100103
$ws.Range("H175").Value = "Frutos de hueso (carozo)"
$ws.Range("I175").Value = 100103006
$ws.Range("J175").Value = "Nectarín"
$ws.Range("K175").Value = "Venus"
$ws.Range("L175").Value = "Primera"
$ws.Range("M175").Value = 270
$ws.Range("N175").Value = 11000
$ws.Range("O175").Value = 12000
$ws.Range("P175").Value = 11444
$ws.Range("Q175").Value = "$/caja 16 kilos empedrada"
$ws.Range("R175").Value = "Región de O'Higgins"
$ws.Range("S175").Value = 715
$ws.Range("T175").Value = 16
